$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (B,D,E,F,G,H,J,K,M) mapped to column indices
$colIndices = @(2, 4, 5, 6, 7, 8, 10, 11, 13)

# New loading-percent values per data row (rows 2-25), one array per row
# in the order of $colIndices above.
$newValues = @{
    2 = @(7.856443890422462, 7.539948643108398, 12.91172326816818, 40.44503893972541, 47.67543407820811, 18.88094980744684, 10.28388871155429, 12.5791833399799, 16.32042164378778)
    3 = @(7.788370333197111, 7.534250498574903, 12.92320591715899, 40.43293257620269, 47.58024016031895, 18.91472696884147, 10.30602710266872, 12.32481071260331, 16.23887285971704)
    4 = @(7.748163425617929, 7.531654661964176, 12.932183569048, 40.43526633481816, 47.53552390375134, 18.93900352049577, 10.3207629537522, 12.16879417544929, 16.19174954099896)
    5 = @(7.732196216424433, 7.530825159981072, 12.93632694527289, 40.43867186751931, 47.52076074962796, 18.94978407603123, 10.32705553062999, 12.10535161984556, 16.17330287596278)
    6 = @(7.729570580844615, 7.530701248131333, 12.93704424349481, 40.43938552432106, 47.51851835324645, 18.95162773587731, 10.32811778529027, 12.09482796115275, 16.17028593417796)
    7 = @(7.747946373563626, 7.531642548882905, 12.93223748437929, 40.43530232810082, 47.53531079233225, 18.9391453192756, 10.32084665286468, 12.16793789452293, 16.19149768061122)
    8 = @(7.832652731690819, 7.537797335891697, 12.9152826810709, 40.43883705393478, 47.63976547377515, 18.89186097938398, 10.29128498436454, 12.49149564477929, 16.2917011758748)
    9 = @(8.010578822739857, 7.556971179373716, 12.89731110931503, 40.52326891393871, 47.95314586734551, 18.82727834298449, 10.24237179803418, 13.12336615495061, 16.51087997724725)
    10 = @(8.147388468622816, 7.575305900899536, 12.89339608473489, 40.63245616131218, 48.24874989096435, 18.7970807500358, 10.21194271131976, 13.5808970477542, 16.6847262736816)
    11 = @(8.21070239706331, 7.584549772283718, 12.89362464766937, 40.69231334081234, 48.39716946127755, 18.78710584823202, 10.19929263910773, 13.78657775134347, 16.76636450645819)
    12 = @(8.234811863677011, 7.588178276354866, 12.89399931701362, 40.71643704082742, 48.4553505630797, 18.78387062076109, 10.1946736036676, 13.864033088462, 16.79762556749236)
    13 = @(8.229613843903611, 7.58739114826392, 12.89390582597771, 40.71117690944565, 48.4427327386183, 18.78454326230905, 10.19566078101488, 13.84737209032145, 16.79087786989924)
    14 = @(8.212683336498197, 7.584845736036601, 12.89364970293303, 40.69426886106538, 48.40191658908085, 18.78682881526847, 10.1989091967375, 13.79295913062488, 16.76892955151684)
    15 = @(8.202329712823209, 7.583303217336744, 12.89353031430311, 40.68410168633289, 48.37717220331567, 18.78829940147699, 10.20092124465854, 13.75957117188508, 16.75553005823253)
    16 = @(8.143270643508561, 7.574719813371759, 12.89342151450974, 40.62874871859596, 48.23932869063096, 18.79780843250922, 10.21279339994054, 13.56739913470093, 16.67944065399792)
    17 = @(8.107300971289479, 7.569684247652584, 12.89386889071879, 40.59739566263209, 48.15831986846549, 18.80460630841523, 10.22038182825857, 13.44882211148172, 16.63340187857533)
    18 = @(8.086714962975703, 7.566873063480998, 12.89431541434327, 40.58032167592975, 48.11304047142769, 18.80887032609533, 10.2248587243549, 13.38039333379507, 16.60716334255458)
    19 = @(8.07976319543627, 7.565935922535163, 12.8944991195881, 40.57470571652987, 48.09793624195085, 18.8103748235834, 10.22639380585148, 13.35718811876466, 16.59832157412886)
    20 = @(8.111119518971831, 7.570211493652162, 12.89380169056673, 40.60063399948461, 48.16680751592882, 18.80384601072535, 10.21956241268992, 13.46146885673014, 16.63827791438689)
    21 = @(8.217652773551464, 7.585589925171295, 12.89371712042302, 40.69919569095469, 48.41385184348413, 18.78614277441528, 10.19795041141714, 13.8089538507108, 16.77536707202201)
    22 = @(8.288048775609369, 7.596386081712712, 12.89534073910011, 40.77209960448631, 48.58682249834733, 18.777732520767, 10.1848239431476, 14.0335072142085, 16.8669719715102)
    23 = @(8.250413636710833, 7.590556390873187, 12.89432087657062, 40.7324157856072, 48.49346149028062, 18.78193180980629, 10.1917385070257, 13.91391681724752, 16.81790391988574)
    24 = @(8.109392860916955, 7.569972864494027, 12.89383148198439, 40.59916698344061, 48.16296621863778, 18.80418863292842, 10.2199325147019, 13.4557520584588, 16.63607274128008)
    25 = @(7.961293427599961, 7.551031833979346, 12.90053949268274, 40.49213358781069, 47.85681595078675, 18.84172632163013, 10.25463592828633, 12.95325025499663, 16.44925927146274)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $colIndices.Length; $i++) {
        $ws.Cells.Item($row, $colIndices[$i]).Value = $rowValues[$i]
    }
}
